$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# Fix wording of the existing MODBUS "button(i,b)" description (row 117, column C):
# "value b" -> "value of b"
$ws.Range("C117").Value = "sets button i to pressed if value of b is yes, true, t, or 1, otherwise to normal"

# Insert a new row for the `visible(i,b)` Artisan Command, right after the
# existing button()/button(<bool>)/button() trio (rows 117-119), pushing
# everything below down by one row.
$ws.Rows.Item(120).Insert()
$ws.Range("B120").Value = "visible(i,b)"
$ws.Range("C120").Value = "sets button i to visible if value of b is yes, true, t, or 1, otherwise to hidden"

# Make the "Commands" sheet the active sheet/tab, with the new row in view.
$ws.Activate()
$ws.Range("C120").Select()
